$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecaster values to insert into a fresh column B (rows 2-16).
$newValues = @{
    2  = -0.5825945370336409
    3  = 0.09567504080935779
    4  = -0.2604190369987228
    5  = 0.8354549961584912
    6  = -0.1000793599026215
    7  = -0.3537865060796963
    8  = 0.1481773904324453
    9  = 0.157445989004155
    10 = -0.5006594565260708
    11 = 0.2803578805354692
    12 = -0.1719748578450117
    13 = 0.3058625397463315
    14 = -0.6123299526872862
    15 = 0.6883713851991116
    16 = -0.2766911554241067
}

# Shift the existing data (columns B:K, rows 2-16) one column to the right,
# dropping any value that would fall past column K, then drop in the new
# naive-forecast column B values.
for ($r = 16; $r -ge 2; $r--) {
    for ($c = 11; $c -ge 3; $c--) {
        $srcCell = $ws.Cells.Item($r, $c - 1)
        $dstCell = $ws.Cells.Item($r, $c)
        $dstCell.Value = $srcCell.Value2
    }
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
